$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '45.172.04'
$ws.Range("E2").Value = '  +3.02%  '
$ws.Range("D3").Value = '2.364.01'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'310.44"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = "'108.15"
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").Value = "'40.75"
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").Value = "'0.0914"
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = "'8.42"
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").Value = "'0.975"
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("D15").Value = '2.723.70'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = "'15.17"
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").Value = '2.364.04'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '45.148.13'
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("D19").Value = "'14.51"
$ws.Range("E19").Value = '  +11.17%  '
$ws.Range("E20").Value = '  -4.74%  '
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = "'72.93"
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").Value = "'3.48"
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").Value = "'259.01"
$ws.Range("E24").Value = '  -3.51%  '
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("E28").Value = '  -6.74%  '
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").Value = "'0.0966"
$ws.Range("E30").Value = '  +8.87%  '
$ws.Range("D31").Value = "'22.30"
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").Value = "'37.04"
$ws.Range("E32").Value = '  -5.64%  '
$ws.Range("D33").Value = "'167.99"
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("E34").Value = '  +5.05%  '
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").Value = "'4.66"
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("D38").Value = "'3.94"
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.0351"
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = "'2.89"
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").Value = "'99.76"
$ws.Range("D43").Value = '1.898.52'
$ws.Range("E43").Value = '  +13.95%  '
$ws.Range("D44").Value = "'69.24"
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("D45").Value = "'0.228"
$ws.Range("E45").Value = '  -4.51%  '
$ws.Range("D46").Value = "'12.80"
$ws.Range("E46").Value = '  -5.01%  '
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").Value = "'81.86"
$ws.Range("E48").Value = '  +6.39%  '
$ws.Range("E49").Value = '  +7.83%  '
$ws.Range("D50").Value = "'110.15"
$ws.Range("E50").Value = '  -3.41%  '
$ws.Range("E51").Value = '  +2.12%  '
